# Daily cryptos-list refresh (GitHub Actions bot).
# Coin/Link/Price/Volume columns are plain text in the sheet, so numeric-looking
# Price values are written with a leading apostrophe to force text (matches the
# workbook's existing inlineStr cells instead of letting Excel re-type them as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.455.40'
$ws.Range("E2").Value = '  +3.45%  '

$ws.Range("D3").Value = '3.072.79'
$ws.Range("E3").Value = '  +5.73%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '''514.57'
$ws.Range("E5").Value = '  +3.39%  '

$ws.Range("D6").Value = '''142.83'
$ws.Range("E6").Value = '  +8.38%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '''0.434'
$ws.Range("E8").Value = '  +3.29%  '

$ws.Range("D9").Value = '''7.28'
$ws.Range("E9").Value = '  +2.27%  '

$ws.Range("E10").Value = '  +4.73%  '

$ws.Range("D11").Value = '''0.372'
$ws.Range("E11").Value = '  +6.77%  '

$ws.Range("D12").Value = '3.601.58'
$ws.Range("E12").Value = '  +5.47%  '

$ws.Range("E13").Value = '  +2.99%  '

$ws.Range("D14").Value = '''26.04'
$ws.Range("E14").Value = '  +1.91%  '

$ws.Range("D15").Value = '''0.0000164'
$ws.Range("E15").Value = '  +4.05%  '

$ws.Range("D16").Value = '57.552.34'
$ws.Range("E16").Value = '  +3.70%  '

$ws.Range("D17").Value = '3.073.15'
$ws.Range("E17").Value = '  +5.69%  '

$ws.Range("D18").Value = '''6.12'
$ws.Range("E18").Value = '  +3.37%  '

$ws.Range("D19").Value = '''12.99'
$ws.Range("E19").Value = '  +3.48%  '

$ws.Range("E20").Value = '  +7.27%  '

$ws.Range("D21").Value = '''336.95'
$ws.Range("E21").Value = '  +8.17%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("D23").Value = '''0.499'
$ws.Range("E23").Value = '  +3.16%  '

$ws.Range("D24").Value = '''65.45'
$ws.Range("E24").Value = '  +4.29%  '

$ws.Range("E25").Value = '  +7.32%  '

$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").Value = '0.0₃0937'
$ws.Range("E27").Value = '  +13.89%  '

$ws.Range("D28").Value = '''6.45'
$ws.Range("E28").Value = '  +3.32%  '

$ws.Range("D29").Value = '''7.13'
$ws.Range("E29").Value = '  +6.02%  '

$ws.Range("D30").Value = '''1.81'
$ws.Range("E30").Value = '  +3.76%  '

$ws.Range("D31").Value = '''20.77'
$ws.Range("E31").Value = '  +5.43%  '

$ws.Range("E32").Value = '  +5.14%  '

$ws.Range("D33").Value = '''154.52'
$ws.Range("E33").Value = '  +1.67%  '

$ws.Range("D34").Value = '''4.54'
$ws.Range("E34").Value = '  +4.88%  '

$ws.Range("D35").Value = '''5.88'
$ws.Range("E35").Value = '  +5.92%  '

$ws.Range("D36").Value = '''25.98'
$ws.Range("E36").Value = '  +10.14%  '

$ws.Range("E37").Value = '  +5.77%  '

$ws.Range("D38").Value = '''0.0679'
$ws.Range("E38").Value = '  +6.34%  '

$ws.Range("D39").Value = '3.110.70'
$ws.Range("E39").Value = '  +5.72%  '

$ws.Range("D40").Value = '''36.87'
$ws.Range("E40").Value = '  +1.22%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '''3.86'
$ws.Range("E41").Value = '  +5.26%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '''0.668'
$ws.Range("E42").Value = '  +5.15%  '

$ws.Range("D43").Value = '''0.999'

$ws.Range("D44").Value = '2.269.89'
$ws.Range("E44").Value = '  +8.26%  '

$ws.Range("E45").Value = '  +5.38%  '

$ws.Range("E46").Value = '  +8.10%  '

$ws.Range("D47").Value = '''0.954'
$ws.Range("E47").Value = '  +4.81%  '

$ws.Range("D48").Value = '''20.36'
$ws.Range("E48").Value = '  +10.41%  '

$ws.Range("D49").Value = '''5.87'
$ws.Range("E49").Value = '  -1.43%  '

$ws.Range("D50").Value = '''0.0873'
$ws.Range("E50").Value = '  +4.52%  '

$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '''0.688'
$ws.Range("E51").Value = '  +7.83%  '
